# Refresh cached Universalis market-board price/profit figures for several
# Leve rows across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets (scheduled runner sync).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86,8).Value = 321434270  # H86
$ws.Cells.Item(86,10).Value = 416669000  # J86
$ws.Cells.Item(86,12).Value = 416669000  # L86
$ws.Cells.Item(86,14).Value = -416671246  # N86
$ws.Cells.Item(89,8).Value = 321434270  # H89
$ws.Cells.Item(89,10).Value = 416669000  # J89
$ws.Cells.Item(89,12).Value = 2083345000  # L89
$ws.Cells.Item(89,14).Value = -2083356232  # N89
$ws.Cells.Item(94,8).Value = 355.83334  # H94
$ws.Cells.Item(94,9).Value = 384.6  # I94
$ws.Cells.Item(94,11).Value = 384.6  # K94
$ws.Cells.Item(94,13).Value = 66.39999999999998  # M94
$ws.Cells.Item(112,8).Value = 3365.453  # H112
$ws.Cells.Item(112,9).Value = 1725  # I112
$ws.Cells.Item(112,10).Value = 3499.3674  # J112
$ws.Cells.Item(112,11).Value = 5175  # K112
$ws.Cells.Item(112,12).Value = 10498.1022  # L112
$ws.Cells.Item(112,13).Value = -4067  # M112
$ws.Cells.Item(112,14).Value = -12714.1022  # N112
$ws.Cells.Item(125,8).Value = 2333.1667  # H125
$ws.Cells.Item(125,9).Value = 2599.8  # I125
$ws.Cells.Item(125,11).Value = 23398.2  # K125
$ws.Cells.Item(125,13).Value = -20938.2  # M125
$ws.Cells.Item(132,8).Value = 6164.6523  # H132
$ws.Cells.Item(132,9).Value = 6567.737  # I132
$ws.Cells.Item(132,11).Value = 19703.211  # K132
$ws.Cells.Item(132,13).Value = -17173.211  # M132
$ws.Cells.Item(138,8).Value = 4140.178  # H138
$ws.Cells.Item(138,9).Value = 4376.737  # I138
$ws.Cells.Item(138,10).Value = 4056.9443  # J138
$ws.Cells.Item(138,11).Value = 13130.211  # K138
$ws.Cells.Item(138,12).Value = 12170.8329  # L138
$ws.Cells.Item(138,13).Value = -7990.210999999999  # M138
$ws.Cells.Item(138,14).Value = -22450.8329  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45,8).Value = 3561.1333  # H45
$ws.Cells.Item(45,9).Value = 3547.36  # I45
$ws.Cells.Item(45,11).Value = 3547.36  # K45
$ws.Cells.Item(45,13).Value = -3170.36  # M45
$ws.Cells.Item(61,8).Value = 4337.5  # H61
$ws.Cells.Item(61,9).Value = 4233.3335  # I61
$ws.Cells.Item(61,11).Value = 4233.3335  # K61
$ws.Cells.Item(61,13).Value = -4021.3335  # M61
$ws.Cells.Item(122,8).Value = 5753.6055  # H122
$ws.Cells.Item(122,9).Value = 4340.96  # I122
$ws.Cells.Item(122,11).Value = 13022.88  # K122
$ws.Cells.Item(122,13).Value = -10572.88  # M122
$ws.Cells.Item(132,8).Value = 154469.95  # H132
$ws.Cells.Item(132,9).Value = 219696.7  # I132
$ws.Cells.Item(132,11).Value = 659090.1000000001  # K132
$ws.Cells.Item(132,13).Value = -656560.1000000001  # M132
$ws.Cells.Item(136,8).Value = 4337.5  # H136
$ws.Cells.Item(136,9).Value = 4233.3335  # I136
$ws.Cells.Item(136,11).Value = 12700.0005  # K136
$ws.Cells.Item(136,13).Value = -10150.0005  # M136
$ws.Cells.Item(138,8).Value = 113904.164  # H138
$ws.Cells.Item(138,9).Value = 99999.25  # I138
$ws.Cells.Item(138,11).Value = 99999.25  # K138
$ws.Cells.Item(138,13).Value = -94859.25  # M138

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134,8).Value = 3761945.2  # H134
$ws.Cells.Item(134,9).Value = 5496551  # I134
$ws.Cells.Item(134,11).Value = 16489653  # K134
$ws.Cells.Item(134,13).Value = -16487118  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22,8).Value = 280.70587  # H22
$ws.Cells.Item(22,9).Value = 368  # I22
$ws.Cells.Item(22,10).Value = 203.11111  # J22
$ws.Cells.Item(22,11).Value = 368  # K22
$ws.Cells.Item(22,12).Value = 203.11111  # L22
$ws.Cells.Item(22,13).Value = -18  # M22
$ws.Cells.Item(22,14).Value = -903.1111100000001  # N22
$ws.Cells.Item(58,8).Value = 3833.6667  # H58
$ws.Cells.Item(58,9).Value = 3664.25  # I58
$ws.Cells.Item(58,10).Value = 4027.2856  # J58
$ws.Cells.Item(58,11).Value = 3664.25  # K58
$ws.Cells.Item(58,12).Value = 4027.2856  # L58
$ws.Cells.Item(58,13).Value = -3461.25  # M58
$ws.Cells.Item(58,14).Value = -4433.2856  # N58
$ws.Cells.Item(86,8).Value = 4131.3335  # H86
$ws.Cells.Item(86,9).Value = 4995  # I86
$ws.Cells.Item(86,10).Value = 3699.5  # J86
$ws.Cells.Item(86,11).Value = 4995  # K86
$ws.Cells.Item(86,12).Value = 3699.5  # L86
$ws.Cells.Item(86,13).Value = -3872  # M86
$ws.Cells.Item(86,14).Value = -5945.5  # N86
$ws.Cells.Item(89,8).Value = 4131.3335  # H89
$ws.Cells.Item(89,9).Value = 4995  # I89
$ws.Cells.Item(89,10).Value = 3699.5  # J89
$ws.Cells.Item(89,11).Value = 24975  # K89
$ws.Cells.Item(89,12).Value = 18497.5  # L89
$ws.Cells.Item(89,13).Value = -19359  # M89
$ws.Cells.Item(89,14).Value = -29729.5  # N89
$ws.Cells.Item(109,8).Value = 0  # H109
$ws.Cells.Item(109,10).Value = 0  # J109
$ws.Cells.Item(109,12).Value = 0  # L109
$ws.Cells.Item(109,14).ClearContents()  # N109
$ws.Cells.Item(134,8).Value = 2381  # H134
$ws.Cells.Item(134,9).Value = 1766.5454  # I134
$ws.Cells.Item(134,10).Value = 3883  # J134
$ws.Cells.Item(134,11).Value = 5299.6362  # K134
$ws.Cells.Item(134,12).Value = 11649  # L134
$ws.Cells.Item(134,13).Value = -2764.6362  # M134
$ws.Cells.Item(134,14).Value = -16719  # N134
$ws.Cells.Item(136,8).Value = 3833.6667  # H136
$ws.Cells.Item(136,9).Value = 3664.25  # I136
$ws.Cells.Item(136,10).Value = 4027.2856  # J136
$ws.Cells.Item(136,11).Value = 10992.75  # K136
$ws.Cells.Item(136,12).Value = 12081.8568  # L136
$ws.Cells.Item(136,13).Value = -8442.75  # M136
$ws.Cells.Item(136,14).Value = -17181.8568  # N136

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(36,8).Value = 4313.077  # H36
$ws.Cells.Item(36,9).Value = 1090.5  # I36
$ws.Cells.Item(36,10).Value = 5745.3335  # J36
$ws.Cells.Item(36,11).Value = 1090.5  # K36
$ws.Cells.Item(36,12).Value = 5745.3335  # L36
$ws.Cells.Item(36,13).Value = -605.5  # M36
$ws.Cells.Item(36,14).Value = -6715.3335  # N36
$ws.Cells.Item(93,8).Value = 49999.5  # H93
$ws.Cells.Item(93,10).Value = 49999.5  # J93
$ws.Cells.Item(93,12).Value = 49999.5  # L93
$ws.Cells.Item(93,14).Value = -53743.5  # N93
$ws.Cells.Item(97,8).Value = 1798.25  # H97
$ws.Cells.Item(97,9).Value = 1369.8572  # I97
$ws.Cells.Item(97,10).Value = 2398  # J97
$ws.Cells.Item(97,11).Value = 1369.8572  # K97
$ws.Cells.Item(97,12).Value = 2398  # L97
$ws.Cells.Item(97,13).Value = -873.8571999999999  # M97
$ws.Cells.Item(97,14).Value = -3390  # N97
$ws.Cells.Item(102,8).Value = 2065.125  # H102
$ws.Cells.Item(102,9).Value = 1788.7858  # I102
$ws.Cells.Item(102,11).Value = 1788.7858  # K102
$ws.Cells.Item(102,13).Value = -166.7858000000001  # M102
$ws.Cells.Item(126,8).Value = 2234.9092  # H126
$ws.Cells.Item(126,9).Value = 2249  # I126
$ws.Cells.Item(126,10).Value = 2231.7778  # J126
$ws.Cells.Item(126,11).Value = 6747  # K126
$ws.Cells.Item(126,12).Value = 6695.3334  # L126
$ws.Cells.Item(126,13).Value = -4277  # M126
$ws.Cells.Item(126,14).Value = -11635.3334  # N126
$ws.Cells.Item(132,8).Value = 2089.3572  # H132
$ws.Cells.Item(132,9).Value = 1936.3684  # I132
$ws.Cells.Item(132,11).Value = 5809.1052  # K132
$ws.Cells.Item(132,13).Value = -3279.1052  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value = 21479  # H7
$ws.Cells.Item(7,9).Value = 2874  # I7
$ws.Cells.Item(7,10).Value = 25200  # J7
$ws.Cells.Item(7,11).Value = 2874  # K7
$ws.Cells.Item(7,12).Value = 25200  # L7
$ws.Cells.Item(7,13).Value = -2762  # M7
$ws.Cells.Item(7,14).Value = -25424  # N7
$ws.Cells.Item(40,8).Value = 18524440  # H40
$ws.Cells.Item(40,10).Value = 5044  # J40
$ws.Cells.Item(40,12).Value = 5044  # L40
$ws.Cells.Item(40,14).Value = -5316  # N40
$ws.Cells.Item(122,8).Value = 2032560.1  # H122
$ws.Cells.Item(122,9).Value = 3835317.5  # I122
$ws.Cells.Item(122,10).Value = 4458.125  # J122
$ws.Cells.Item(122,11).Value = 11505952.5  # K122
$ws.Cells.Item(122,12).Value = 13374.375  # L122
$ws.Cells.Item(122,13).Value = -11503502.5  # M122
$ws.Cells.Item(122,14).Value = -18274.375  # N122
$ws.Cells.Item(126,8).Value = 21479  # H126
$ws.Cells.Item(126,9).Value = 2874  # I126
$ws.Cells.Item(126,10).Value = 25200  # J126
$ws.Cells.Item(126,11).Value = 8622  # K126
$ws.Cells.Item(126,12).Value = 75600  # L126
$ws.Cells.Item(126,13).Value = -6152  # M126
$ws.Cells.Item(126,14).Value = -80540  # N126
$ws.Cells.Item(132,8).Value = 155863.66  # H132
$ws.Cells.Item(132,9).Value = 592300.75  # I132
$ws.Cells.Item(132,11).Value = 1776902.25  # K132
$ws.Cells.Item(132,13).Value = -1774372.25  # M132
$ws.Cells.Item(136,8).Value = 5388.643  # H136
$ws.Cells.Item(136,9).Value = 5120.125  # I136
$ws.Cells.Item(136,10).Value = 6999.75  # J136
$ws.Cells.Item(136,11).Value = 15360.375  # K136
$ws.Cells.Item(136,12).Value = 20999.25  # L136
$ws.Cells.Item(136,13).Value = -12810.375  # M136
$ws.Cells.Item(136,14).Value = -26099.25  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(109,8).Value = 79971.75  # H109
$ws.Cells.Item(109,10).Value = 79971.75  # J109
$ws.Cells.Item(109,12).Value = 79971.75  # L109
$ws.Cells.Item(109,14).Value = -82745.75  # N109
$ws.Cells.Item(122,8).Value = 142860130  # H122
$ws.Cells.Item(122,10).Value = 4200  # J122
$ws.Cells.Item(122,12).Value = 12600  # L122
$ws.Cells.Item(122,14).Value = -17500  # N122
$ws.Cells.Item(126,8).Value = 5130.909  # H126
$ws.Cells.Item(126,9).Value = 2560.7144  # I126
$ws.Cells.Item(126,11).Value = 7682.1432  # K126
$ws.Cells.Item(126,13).Value = -5212.1432  # M126
$ws.Cells.Item(132,8).Value = 38282.9  # H132
$ws.Cells.Item(132,9).Value = 44918.918  # I132
$ws.Cells.Item(132,11).Value = 134756.754  # K132
$ws.Cells.Item(132,13).Value = -132226.754  # M132
$ws.Cells.Item(136,8).Value = 44909.883  # H136
$ws.Cells.Item(136,9).Value = 2507.7144  # I136
$ws.Cells.Item(136,11).Value = 7523.1432  # K136
$ws.Cells.Item(136,13).Value = -4973.1432  # M136
